# Apply 'Updated cryptos list' data refresh to Sheet1
# Values are forced to Text format (quotePrefix-free, style-neutral)
# so numeric-looking strings (e.g. '595.54') do not get coerced to
# real numbers by Excel's automatic type detection, matching the
# original workbook where every data cell is an inline/shared string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '68.130.71'
Set-TextValue 'E2' '  -2.29%  '

# Row 3
Set-TextValue 'D3' '3.796.81'
Set-TextValue 'E3' '  +2.44%  '

# Row 4
Set-TextValue 'E4' '  +0.26%  '

# Row 5
Set-TextValue 'D5' '595.54'
Set-TextValue 'E5' '  -3.48%  '

# Row 6
Set-TextValue 'E6' '  -4.80%  '

# Row 7
Set-TextValue 'D7' '3.794.52'
Set-TextValue 'E7' '  +2.29%  '

# Row 8
Set-TextValue 'E8' '  +0.15%  '

# Row 9
Set-TextValue 'D9' '0.534'
Set-TextValue 'E9' '  +0.60%  '

# Row 10
Set-TextValue 'D10' '0.159'
Set-TextValue 'E10' '  -3.16%  '

# Row 11
Set-TextValue 'D11' '6.32'
Set-TextValue 'E11' '  +0.48%  '

# Row 12
Set-TextValue 'D12' '0.469'
Set-TextValue 'E12' '  -2.37%  '

# Row 13
Set-TextValue 'D13' '38.37'
Set-TextValue 'E13' '  -4.10%  '

# Row 14
Set-TextValue 'D14' '0.0000244'
Set-TextValue 'E14' '  -3.68%  '

# Row 15
Set-TextValue 'D15' '4.433.00'
Set-TextValue 'E15' '  +2.58%  '

# Row 16
Set-TextValue 'D16' '3.796.72'
Set-TextValue 'E16' '  +2.84%  '

# Row 17
Set-TextValue 'D17' '68.291.09'
Set-TextValue 'E17' '  -1.95%  '

# Row 18
Set-TextValue 'B18' 'TRON'
Set-TextValue 'C18' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 'D18' '0.116'
Set-TextValue 'E18' '  -4.42%  '

# Row 19
Set-TextValue 'B19' 'Polkadot'
Set-TextValue 'C19' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D19' '7.26'
Set-TextValue 'E19' '  -4.05%  '

# Row 20
Set-TextValue 'D20' '16.13'
Set-TextValue 'E20' '  -1.37%  '

# Row 21
Set-TextValue 'D21' '488.96'
Set-TextValue 'E21' '  -2.33%  '

# Row 22
Set-TextValue 'D22' '9.39'
Set-TextValue 'E22' '  +2.29%  '

# Row 23
Set-TextValue 'D23' '0.737'
Set-TextValue 'E23' '  +1.57%  '

# Row 24
Set-TextValue 'D24' '86.12'
Set-TextValue 'E24' '  -0.06%  '

# Row 25
Set-TextValue 'D25' '2.38'
Set-TextValue 'E25' '  -5.14%  '

# Row 26
Set-TextValue 'E26' '  +6.42%  '

# Row 27
Set-TextValue 'D27' '12.27'
Set-TextValue 'E27' '  -5.42%  '

# Row 28
Set-TextValue 'D28' '10.25'
Set-TextValue 'E28' '  -7.73%  '

# Row 29
Set-TextValue 'D29' '1.00'
Set-TextValue 'E29' '  -0.22%  '

# Row 30
Set-TextValue 'D30' '2.94'
Set-TextValue 'E30' '  +1.40%  '

# Row 31
Set-TextValue 'E31' '  -0.32%  '

# Row 32
Set-TextValue 'D32' '32.51'
Set-TextValue 'E32' '  +7.51%  '

# Row 33
Set-TextValue 'D33' '7.62'
Set-TextValue 'E33' '  -4.46%  '

# Row 34
Set-TextValue 'E34' '  -2.83%  '

# Row 35
Set-TextValue 'D35' '0.999'
Set-TextValue 'E35' '  +0.07%  '

# Row 36
Set-TextValue 'E36' '  -2.41%  '

# Row 37
Set-TextValue 'D37' '5.85'
Set-TextValue 'E37' '  -3.30%  '

# Row 38
Set-TextValue 'E38' '  -1.62%  '

# Row 39
Set-TextValue 'D39' '0.325'
Set-TextValue 'E39' '  -4.47%  '

# Row 40
Set-TextValue 'D40' '451.04'
Set-TextValue 'E40' '  +5.53%  '

# Row 41
Set-TextValue 'D41' '49.14'
Set-TextValue 'E41' '  -1.78%  '

# Row 42
Set-TextValue 'E42' '  -2.22%  '

# Row 43
Set-TextValue 'E43' '  -2.96%  '

# Row 44
Set-TextValue 'D44' '8.35'
Set-TextValue 'E44' '  -2.77%  '

# Row 45
Set-TextValue 'D45' '41.46'
Set-TextValue 'E45' '  -5.58%  '

# Row 46
Set-TextValue 'D46' '2.859.06'
Set-TextValue 'E46' '  -2.78%  '

# Row 47
Set-TextValue 'B47' 'VeChain'
Set-TextValue 'C47' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D47' '0.0353'
Set-TextValue 'E47' '  -1.62%  '

# Row 48
Set-TextValue 'B48' 'USDe'
Set-TextValue 'C48' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D48' '1.00'
Set-TextValue 'E48' '  +0.08%  '

# Row 49
Set-TextValue 'D49' '137.37'
Set-TextValue 'E49' '  +0.43%  '

# Row 50
Set-TextValue 'D50' '26.56'
Set-TextValue 'E50' '  -2.87%  '

# Row 51
Set-TextValue 'E51' '  +8.57%  '

